# Auto-generated edit script: update TPM-derived NATMI metrics for Efna1-Epha7 sheet
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 23.73148533333334
$ws.Range("H2").Value = 71.194456
$ws.Range("I2").Value = 0.8653076146801144
$ws.Range("J2").Value = 0.8653076146801145
$ws.Range("M2").Value = 0.1352566666666667
$ws.Range("N2").Value = 0.40577
$ws.Range("O2").Value = 0.1173241749329269
$ws.Range("P2").Value = 0.1173241749329268
$ws.Range("Q2").Value = 3.209841601235556
$ws.Range("R2").Value = 28.88857441112
$ws.Range("S2").Value = 0.1015215019555234
$ws.Range("T2").Value = 0.1015215019555234
$ws.Range("G3").Value = 23.73148533333334
$ws.Range("H3").Value = 71.194456
$ws.Range("I3").Value = 0.8653076146801144
$ws.Range("J3").Value = 0.8653076146801145
$ws.Range("O3").Value = 0.03951584152489912
$ws.Range("P3").Value = 0.03951584152489912
$ws.Range("Q3").Value = 1.081103635350222
$ws.Range("R3").Value = 9.729932718152002
$ws.Range("S3").Value = 0.03419335857198787
$ws.Range("T3").Value = 0.03419335857198787
$ws.Range("G4").Value = 23.73148533333334
$ws.Range("H4").Value = 71.194456
$ws.Range("I4").Value = 0.8653076146801144
$ws.Range("J4").Value = 0.8653076146801145
$ws.Range("K4").Value = 3
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 0.9440163333333332
$ws.Range("N4").Value = 2.832049
$ws.Range("O4").Value = 0.8188575111383802
$ws.Range("P4").Value = 0.8188575111383801
$ws.Range("Q4").Value = 22.40290976892711
$ws.Range("R4").Value = 201.626187920344
$ws.Range("S4").Value = 0.708563639726047
$ws.Range("T4").Value = 0.708563639726047
$ws.Range("G5").Value = 23.73148533333334
$ws.Range("H5").Value = 71.194456
$ws.Range("I5").Value = 0.8653076146801144
$ws.Range("J5").Value = 0.8653076146801145
$ws.Range("K5").Value = 2
$ws.Range("L5").Value = 0.6666666666666666
$ws.Range("M5").Value = 0.028017
$ws.Range("N5").Value = 0.084051
$ws.Range("O5").Value = 0.02430247240379386
$ws.Range("P5").Value = 0.02430247240379386
$ws.Range("Q5").Value = 0.664885024584
$ws.Range("R5").Value = 5.983965221256001
$ws.Range("S5").Value = 0.02102911442655617
$ws.Range("T5").Value = 0.02102911442655617
$ws.Range("I6").Value = 0.09010639372350319
$ws.Range("J6").Value = 0.09010639372350321
$ws.Range("M6").Value = 0.1352566666666667
$ws.Range("N6").Value = 0.40577
$ws.Range("O6").Value = 0.1173241749329269
$ws.Range("P6").Value = 0.1173241749329268
$ws.Range("Q6").Value = 0.3342478977466666
$ws.Range("R6").Value = 3.00823107972
$ws.Range("S6").Value = 0.01057165829979147
$ws.Range("T6").Value = 0.01057165829979147
$ws.Range("I7").Value = 0.09010639372350319
$ws.Range("J7").Value = 0.09010639372350321
$ws.Range("O7").Value = 0.03951584152489912
$ws.Range("P7").Value = 0.03951584152489912
$ws.Range("S7").Value = 0.003560629974758117
$ws.Range("T7").Value = 0.003560629974758118
$ws.Range("I8").Value = 0.09010639372350319
$ws.Range("J8").Value = 0.09010639372350321
$ws.Range("K8").Value = 3
$ws.Range("L8").Value = 1
$ws.Range("M8").Value = 0.9440163333333332
$ws.Range("N8").Value = 2.832049
$ws.Range("O8").Value = 0.8188575111383802
$ws.Range("P8").Value = 0.8188575111383801
$ws.Range("Q8").Value = 2.332864491129333
$ws.Range("R8").Value = 20.995780420164
$ws.Range("S8").Value = 0.07378429730208279
$ws.Range("T8").Value = 0.07378429730208279
$ws.Range("I9").Value = 0.09010639372350319
$ws.Range("J9").Value = 0.09010639372350321
$ws.Range("K9").Value = 2
$ws.Range("L9").Value = 0.6666666666666666
$ws.Range("M9").Value = 0.028017
$ws.Range("N9").Value = 0.084051
$ws.Range("O9").Value = 0.02430247240379386
$ws.Range("P9").Value = 0.02430247240379386
$ws.Range("Q9").Value = 0.069235946604
$ws.Range("R9").Value = 0.6231235194360001
$ws.Range("S9").Value = 0.002189808146870821
$ws.Range("T9").Value = 0.002189808146870821
$ws.Range("G10").Value = 1.148663
$ws.Range("H10").Value = 3.445989
$ws.Range("I10").Value = 0.04188304383987305
$ws.Range("J10").Value = 0.04188304383987305
$ws.Range("M10").Value = 0.1352566666666667
$ws.Range("N10").Value = 0.40577
$ws.Range("O10").Value = 0.1173241749329269
$ws.Range("P10").Value = 0.1173241749329268
$ws.Range("Q10").Value = 0.1553643285033333
$ws.Range("R10").Value = 1.39827895653
$ws.Range("S10").Value = 0.004913893562192711
$ws.Range("T10").Value = 0.004913893562192709
$ws.Range("G11").Value = 1.148663
$ws.Range("H11").Value = 3.445989
$ws.Range("I11").Value = 0.04188304383987305
$ws.Range("J11").Value = 0.04188304383987305
$ws.Range("O11").Value = 0.03951584152489912
$ws.Range("P11").Value = 0.03951584152489912
$ws.Range("Q11").Value = 0.05232810874033333
$ws.Range("R11").Value = 0.470952978663
$ws.Range("S11").Value = 0.001655043722956826
$ws.Range("T11").Value = 0.001655043722956826
$ws.Range("G12").Value = 1.148663
$ws.Range("H12").Value = 3.445989
$ws.Range("I12").Value = 0.04188304383987305
$ws.Range("J12").Value = 0.04188304383987305
$ws.Range("K12").Value = 3
$ws.Range("L12").Value = 1
$ws.Range("M12").Value = 0.9440163333333332
$ws.Range("N12").Value = 2.832049
$ws.Range("O12").Value = 0.8188575111383802
$ws.Range("P12").Value = 0.8188575111383801
$ws.Range("Q12").Value = 1.084356633495666
$ws.Range("R12").Value = 9.759209701460998
$ws.Range("S12").Value = 0.03429624503761811
$ws.Range("T12").Value = 0.03429624503761811
$ws.Range("G13").Value = 1.148663
$ws.Range("H13").Value = 3.445989
$ws.Range("I13").Value = 0.04188304383987305
$ws.Range("J13").Value = 0.04188304383987305
$ws.Range("K13").Value = 2
$ws.Range("L13").Value = 0.6666666666666666
$ws.Range("M13").Value = 0.028017
$ws.Range("N13").Value = 0.084051
$ws.Range("O13").Value = 0.02430247240379386
$ws.Range("P13").Value = 0.02430247240379386
$ws.Range("Q13").Value = 0.032182091271
$ws.Range("R13").Value = 0.289638821439
$ws.Range("S13").Value = 0.001017861517105403
$ws.Range("T13").Value = 0.001017861517105403
$ws.Range("G14").Value = 0.07412966666666666
$ws.Range("H14").Value = 0.222389
$ws.Range("I14").Value = 0.002702947756509242
$ws.Range("J14").Value = 0.002702947756509243
$ws.Range("M14").Value = 0.1352566666666667
$ws.Range("N14").Value = 0.40577
$ws.Range("O14").Value = 0.1173241749329269
$ws.Range("P14").Value = 0.1173241749329268
$ws.Range("Q14").Value = 0.01002653161444444
$ws.Range("R14").Value = 0.09023878452999999
$ws.Range("S14").Value = 0.0003171211154192526
$ws.Range("T14").Value = 0.0003171211154192525
$ws.Range("G15").Value = 0.07412966666666666
$ws.Range("H15").Value = 0.222389
$ws.Range("I15").Value = 0.002702947756509242
$ws.Range("J15").Value = 0.002702947756509243
$ws.Range("O15").Value = 0.03951584152489912
$ws.Range("P15").Value = 0.03951584152489912
$ws.Range("Q15").Value = 0.003377026384777778
$ws.Range("R15").Value = 0.030393237463
$ws.Range("S15").Value = 0.0001068092551963008
$ws.Range("T15").Value = 0.0001068092551963008
$ws.Range("G16").Value = 0.07412966666666666
$ws.Range("H16").Value = 0.222389
$ws.Range("I16").Value = 0.002702947756509242
$ws.Range("J16").Value = 0.002702947756509243
$ws.Range("K16").Value = 3
$ws.Range("L16").Value = 1
$ws.Range("M16").Value = 0.9440163333333332
$ws.Range("N16").Value = 2.832049
$ws.Range("O16").Value = 0.8188575111383802
$ws.Range("P16").Value = 0.8188575111383801
$ws.Range("Q16").Value = 0.06997961611788887
$ws.Range("R16").Value = 0.6298165450609999
$ws.Range("S16").Value = 0.002213329072632227
$ws.Range("T16").Value = 0.002213329072632227
$ws.Range("G17").Value = 0.07412966666666666
$ws.Range("H17").Value = 0.222389
$ws.Range("I17").Value = 0.002702947756509242
$ws.Range("J17").Value = 0.002702947756509243
$ws.Range("K17").Value = 2
$ws.Range("L17").Value = 0.6666666666666666
$ws.Range("M17").Value = 0.028017
$ws.Range("N17").Value = 0.084051
$ws.Range("O17").Value = 0.02430247240379386
$ws.Range("P17").Value = 0.02430247240379386
$ws.Range("Q17").Value = 0.002076890871
$ws.Range("R17").Value = 0.018692017839
$ws.Range("S17").Value = 0.0000656883132614624
$ws.Range("T17").Value = 0.0000656883132614624
